# [All Hosts] (overview) Reorienting overview for the M365 ecosystem
#
# Updates the "Office programmability" overview diagram:
#  - Re-caches the datetimeFigureOut date field text (slide master + all
#    11 slide layouts) from 3/10/2023 to 12/8/2025.
#  - Shrinks/repositions + relabels "End-user approachable" -> "End-user skills"
#  - Relabels "Desktop and individual" -> "Windows only"
#  - Collapses the 4-run "Cross-platform and collaboration" label down to a
#    single run reading "Cross-platform"
#  - Shrinks/repositions + relabels "Developer audience " -> "Developer skills"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer/date placeholders on the slide master and every slide layout
#    (the cached text PowerPoint shows for the "Update automatically"
#    date field) - bump the stale cached date forward.
# ---------------------------------------------------------------------
$newDate = "12/8/2025"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $sh = $layout.Shapes.Item($si)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2) Diagram textbox edits on slide 1
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

# "TextBox 24" - End-user approachable -> End-user skills (also shrinks
# and shifts right now that the copy is shorter).
$endUser = $s.Shapes.Item(3)
$endUser.Left = 444.6647244094488
$endUser.Width = 81.6648131496063
$endUser.TextFrame.TextRange.Text = "End-user skills"

# "TextBox 30" - Desktop and individual -> Windows only
$desktop = $s.Shapes.Item(6)
$desktop.TextFrame.TextRange.Text = "Windows only"

# "TextBox 32" - Cross-platform and collaboration -> Cross-platform
$crossPlatform = $s.Shapes.Item(7)
$crossPlatform.TextFrame.TextRange.Text = "Cross-platform"

# "TextBox 34" - Developer audience -> Developer skills (also shrinks and
# shifts right now that the copy is shorter).
$developer = $s.Shapes.Item(8)
$developer.Left = 440.6029233858267
$developer.Width = 89.7884351968504
$developer.TextFrame.TextRange.Text = "Developer skills"
